$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ADM001 / Aarav Kumar Sharma record, updated message + new "Sent By"
$ws.Range("D2").Value = "04 February 2025, 12:00 AM"
$ws.Range("E2").Value = "Hi Student 1, This is a Test Message"
$ws.Range("G2").Value = "John Smith"

# Row 3: ADM002 / Vivaan Raj Gupta record, updated message + new "Sent By"
$ws.Range("D3").Value = "04 February 2025, 12:00 AM"
$ws.Range("E3").Value = "Hi Student 2, This is a Test Message"
$ws.Range("G3").Value = "John Smith"

# Delete old rows 4 and 5 (ADM003 entry and the duplicate ADM001 "test message" row)
$ws.Range("A4:G5").Delete()
